$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 10
$ws.Range("Q2").Value = 2.07
$ws.Range("R2").Value = 1.83
